$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hub")

# Add the new row 7 values (HubMenu navigation test case)
$ws.Range("A7").Value = "HubMenu"
$ws.Range("J7").Value = 365
$ws.Range("K7").Value = 370

# Update the selected cell on the sheet to C6
$ws.Activate()
$ws.Range("C6").Select()
